$d = $word.ActiveDocument

# Two Pearson-logo pictures live in the footers ("image1.png" -> "image2.png")
# and two BTec-logo pictures live in the headers ("image2.jpg" -> "image1.jpg").
# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2, wdHeaderFooterEvenPages = 3.
# NOTE: the loop variable inside the helper must not share a name with the
# caller's loop variable (this interpreter does not scope them separately).

function Rename-LogoInStory($hf, $newName) {
    for ($k = 1; $k -le $hf.Range.Paragraphs.Count; $k++) {
        if ($hf.Range.Paragraphs($k).Range.InlineShapes.Count -gt 0) {
            $hf.Range.Paragraphs($k).Range.InlineShapes(1).Name = $newName
        }
    }
}

$section = $d.Sections(1)

for ($h = 1; $h -le $section.Headers.Count; $h++) {
    if ($section.Headers($h).Exists -and $section.Headers($h).Range.InlineShapes.Count -gt 0) {
        Rename-LogoInStory $section.Headers($h) "image1.jpg"
    }
}

for ($f = 1; $f -le $section.Footers.Count; $f++) {
    if ($section.Footers($f).Exists -and $section.Footers($f).Range.InlineShapes.Count -gt 0) {
        Rename-LogoInStory $section.Footers($f) "image2.png"
    }
}
